$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").ClearContents()
$ws.Range("D2").Value = 2172
$ws.Range("E2").Value = 135
$ws.Range("F2").Value = 135
$ws.Range("G2").Value = 226
$ws.Range("H2").Value = 175
$ws.Range("I2").Value = 169
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 4294
$ws.Range("L2").Value = 541
$ws.Range("M2").Value = 3753
$ws.Range("N2").Value = 3753
$ws.Range("P2").Value = 165
$ws.Range("Q2").Value = 101
$ws.Range("R2").Value = 130
$ws.Range("S2").Value = -345
$ws.Range("T2").Value = 233
$ws.Range("U2").Value = -132
$ws.Range("V2").Value = 11
$ws.Range("W2").Value = 6.24
$ws.Range("X2").Value = 8.06
$ws.Range("Y2").Value = 4.58
$ws.Range("Z2").Value = 3.95
$ws.Range("AA2").Value = 14.41
$ws.Range("AB2").Value = 2238.95
$ws.Range("AC2").Value = 513
$ws.Range("AD2").Value = 24.95
$ws.Range("AE2").Value = 16485
$ws.Range("AF2").Value = 0.78
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 1.95
$ws.Range("AI2").Value = 33.76
$ws.Range("AJ2").Value = 33000000

# Row 3
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("D3").Value = 2033
$ws.Range("E3").Value = 63
$ws.Range("F3").Value = 63
$ws.Range("G3").Value = 114
$ws.Range("H3").Value = 92
$ws.Range("I3").Value = 92
$ws.Range("K3").Value = 4293
$ws.Range("L3").Value = 497
$ws.Range("M3").Value = 3796
$ws.Range("N3").Value = 3796
$ws.Range("P3").Value = 165
$ws.Range("Q3").Value = 208
$ws.Range("R3").Value = -534
$ws.Range("S3").Value = -63
$ws.Range("T3").Value = 360
$ws.Range("U3").Value = -153
$ws.Range("V3").Value = 5
$ws.Range("W3").Value = 3.09
$ws.Range("X3").Value = 4.54
$ws.Range("Y3").Value = 2.44
$ws.Range("Z3").Value = 2.15
$ws.Range("AA3").Value = 13.08
$ws.Range("AB3").Value = 2263.25
$ws.Range("AC3").Value = 279
$ws.Range("AD3").Value = 38.65
$ws.Range("AE3").Value = 16606
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 280
$ws.Range("AH3").Value = 2.59
$ws.Range("AI3").Value = 69.42
$ws.Range("AJ3").Value = 33000000

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("D4").Value = 1953
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 24
$ws.Range("G4").Value = 210
$ws.Range("H4").Value = 217
$ws.Range("I4").Value = 217
$ws.Range("K4").Value = 4331
$ws.Range("L4").Value = 414
$ws.Range("M4").Value = 3917
$ws.Range("N4").Value = 3917
$ws.Range("P4").Value = 165
$ws.Range("Q4").Value = 147
$ws.Range("R4").Value = -62
$ws.Range("S4").Value = -103
$ws.Range("T4").Value = 121
$ws.Range("U4").Value = 25
$ws.Range("W4").Value = 1.21
$ws.Range("X4").Value = 11.14
$ws.Range("Y4").Value = 5.64
$ws.Range("Z4").Value = 5.04
$ws.Range("AA4").Value = 10.57
$ws.Range("AB4").Value = 2359.62
$ws.Range("AC4").Value = 659
$ws.Range("AD4").Value = 14.04
$ws.Range("AE4").Value = 17418
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 2.7
$ws.Range("AI4").Value = 25.85
$ws.Range("AJ4").Value = 33000000

# Row 5
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("D5").Value = 2115
$ws.Range("E5").Value = 128
$ws.Range("F5").Value = 128
$ws.Range("G5").Value = 172
$ws.Range("H5").Value = 127
$ws.Range("I5").Value = 127
$ws.Range("K5").Value = 4315
$ws.Range("L5").Value = 358
$ws.Range("M5").Value = 3957
$ws.Range("N5").Value = 3957
$ws.Range("P5").Value = 165
$ws.Range("Q5").Value = 189
$ws.Range("R5").Value = -183
$ws.Range("S5").Value = -81
$ws.Range("T5").Value = 35
$ws.Range("U5").Value = 155
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 6.03
$ws.Range("X5").Value = 6
$ws.Range("Y5").Value = 3.23
$ws.Range("Z5").Value = 2.94
$ws.Range("AA5").Value = 9.06
$ws.Range("AB5").Value = 2401.32
$ws.Range("AC5").Value = 385
$ws.Range("AD5").Value = 23.91
$ws.Range("AE5").Value = 17814
$ws.Range("AF5").Value = 0.52
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 3.8
$ws.Range("AI5").Value = 61.23
$ws.Range("AJ5").Value = 33000000

# Row 6
$ws.Range("D6").Value = 1953
$ws.Range("E6").Value = 72
$ws.Range("F6").Value = 72
$ws.Range("G6").Value = 127
$ws.Range("H6").Value = 98
$ws.Range("I6").Value = 98
$ws.Range("K6").Value = 4342
$ws.Range("L6").Value = 390
$ws.Range("M6").Value = 3952
$ws.Range("N6").Value = 3952
$ws.Range("P6").Value = 165
$ws.Range("Q6").Value = 218
$ws.Range("R6").Value = -16
$ws.Range("S6").Value = -89
$ws.Range("T6").Value = 45
$ws.Range("U6").Value = 173
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 3.71
$ws.Range("X6").Value = 5.01
$ws.Range("Y6").Value = 2.48
$ws.Range("Z6").Value = 2.26
$ws.Range("AA6").Value = 9.869999999999999
$ws.Range("AB6").Value = 2408.94
$ws.Range("AC6").Value = 297
$ws.Range("AD6").Value = 33.68
$ws.Range("AE6").Value = 17893
$ws.Range("AF6").Value = 0.5600000000000001
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 3
$ws.Range("AI6").Value = 67.69
$ws.Range("AJ6").Value = 33000000

# Rows 7-9: all financial data cells (D:AJ) are removed, leaving only A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
